$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New SamplesTab query text: drops the Tumor / Analyte Type columns from the
# previous SELECT, keeping the ORDER BY / LIMIT 100 clause.
$newSamplesQuery = @"
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
    s.phs_accession = 'phs001437' AND sp.gender = 'Unknown'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
"@

$ws.Range("B3").Value = $newSamplesQuery

# The TSV/Web file-name cells for the SamplesTab and FilesTab rows are no
# longer populated - clear them out entirely.
$ws.Range("D3:E3").ClearContents()
$ws.Range("D4:E4").ClearContents()

# Move the active selection/top-left cell to reflect the new layout.
$ws.Activate()
$ws.Range("B3").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
